# Data clean-up: add a "No of Recipes" count/header row at the very top of
# the sheet (row 1 was previously blank - the recipe data starts on row 2 -
# so we populate it directly instead of inserting a row, which keeps the
# existing rows 2-11 exactly where they are).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Label for the new row
$ws.Range("A1").Value = "No of Recipes"

# Sequential recipe numbers 1..223, one per recipe column (C, E, G, ... QE),
# matching the same "every other column" layout already used by the
# existing recipe rows below.
$col = 3
for ($i = 1; $i -le 223; $i++) {
    $ws.Cells.Item(1, $col).Value = $i
    $col = $col + 2
}

# Leave the selection where the workbook was last saved from.
$ws.Range("QU1").Select()
